$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("storageAssets")

$ws.Range("J3").Value = 50
$ws.Range("L3").Value = 10000000

$ws.Range("J4").Value = 60
$ws.Range("L4").Value = 10000000

$ws.Range("J5").Value = 70
$ws.Range("L5").Value = 10000000

$ws.Range("J6").Value = 80
$ws.Range("L6").Value = 10000000

$ws.Range("J7").Value = 90
$ws.Range("L7").Value = 10000000

$ws.Range("J8").Value = 100
$ws.Range("L8").Value = 10000000

$ws.Range("J9").Value = 110
$ws.Range("L9").Value = 10000000

$ws.Range("F11").Value = 10000
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 100000000
